# semana 37 de 2025
# Adds the week-37 (column AN) counts, a late week-36 (AM) figure for a few
# facilities, and corrects several previously-entered weekly counts in row 28
# (Centro de Salud Boston) and row 35 (week 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Corrections to previously reported values ---
$ws.Range("T28").Value = 98
$ws.Range("X28").Value = 184
$ws.Range("Y28").Value = 186
$ws.Range("Z28").Value = 205
$ws.Range("AA28").Value = 228
$ws.Range("AB28").Value = 235
$ws.Range("AC28").Value = 124
$ws.Range("AD28").Value = 232
$ws.Range("AE28").Value = 206
$ws.Range("AF28").Value = 231
$ws.Range("AG28").Value = 240
$ws.Range("AH28").Value = 230
$ws.Range("AI28").Value = 224
$ws.Range("AJ28").Value = 211
$ws.Range("AK28").Value = 189
$ws.Range("AL28").Value = 134
$ws.Range("M35").Value = 39

# --- New week 37 column (AN) header ---
# Entered with a leading apostrophe so it is stored as text (matching the
# other week-number headers in row 1), then the formatting is copied from the
# neighbouring header cell so the new header looks identical to the rest.
$ws.Range("AN1").Value = "'37"
$ws.Range("AM1").Copy() | Out-Null
$ws.Range("AN1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- New week 37 (column AN) weekly counts, plus a few late week 36 (AM) values ---
$ws.Range("AN2").Value = 63
$ws.Range("AN3").Value = 72
$ws.Range("AN5").Value = 1
$ws.Range("AN6").Value = 145
$ws.Range("AN7").Value = 28
$ws.Range("AN8").Value = 24
$ws.Range("AN9").Value = 3
$ws.Range("AN10").Value = 6
$ws.Range("AN12").Value = 2
$ws.Range("AN14").Value = 3
$ws.Range("AN15").Value = 3
$ws.Range("AN16").Value = 1
$ws.Range("AN17").Value = 1
$ws.Range("AN21").Value = 1
$ws.Range("AN23").Value = 6
$ws.Range("AN25").Value = 60
$ws.Range("AN26").Value = 1
$ws.Range("AN28").Value = 234
$ws.Range("AN29").Value = 0
$ws.Range("AN30").Value = 25
$ws.Range("AM31").Value = 5
$ws.Range("AN31").Value = 3
$ws.Range("AM35").Value = 41
$ws.Range("AN35").Value = 33
$ws.Range("AN36").Value = 5
$ws.Range("AN37").Value = 4
$ws.Range("AN38").Value = 92
$ws.Range("AN41").Value = 8
$ws.Range("AN42").Value = 14
$ws.Range("AM43").Value = 27
$ws.Range("AN43").Value = 33
$ws.Range("AN44").Value = 12
$ws.Range("AN45").Value = 90
$ws.Range("AN46").Value = 125
$ws.Range("AN47").Value = 6
$ws.Range("AN48").Value = 117
$ws.Range("AN49").Value = 3
$ws.Range("AN50").Value = 0
$ws.Range("AN51").Value = 4
$ws.Range("AM52").Value = 5
$ws.Range("AN53").Value = 14
$ws.Range("AN54").Value = 1
$ws.Range("AN55").Value = 0
$ws.Range("AN56").Value = 2
$ws.Range("AN57").Value = 100
$ws.Range("AN58").Value = 13

Write-Output "Applied semana 37 de 2025 updates"
